$wb = $excel.ActiveWorkbook

# --- 1) "Games" sheet: append the completed game (was the next scheduled
#        game vs NYK on 45306) as row 41 with its final boxscore stats. ---
$games = $wb.Worksheets.Item("Games")

$games.Range("A41").Value = 40
$games.Range("B41").NumberFormat = "YYYY-MM-DD"
$games.Range("B41").Value = 45306
$games.Range("C41").Value = 1
$games.Range("D41").Value = 98
$games.Range("E41").Value = 85.09999999999999
$games.Range("F41").Value = 0.494
$games.Range("G41").Value = 6.9
$games.Range("H41").Value = 31.7
$games.Range("I41").Value = 0.165
$games.Range("J41").Value = 115.2
$games.Range("K41").Value = "NYK"
$games.Range("L41").Value = 94
$games.Range("M41").Value = 0.476
$games.Range("N41").Value = 8.9
$games.Range("O41").Value = 28.3
$games.Range("P41").Value = 0.181
$games.Range("Q41").Value = 110.5
$games.Range("R41").Value = 0
$games.Range("S41").Value = 1

# --- 2) "Next" sheet: that game is no longer upcoming, so remove its row
#        and shift the remaining schedule up by one. ---
$next = $wb.Worksheets.Item("Next")
$next.Rows(2).Delete()
